$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B18 (semana 17) with new case count
$ws.Range("B18").Value = 547

# Add new row for semana 18
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 5
